# Denmark "Bank (Money Center)" capital-structure database refresh.
# The underlying data provider reshuffled which company occupies which
# spreadsheet row and refreshed every financial metric column (D..AM)
# for rows 2-8. Cells that no longer have a reported value are cleared
# (ClearContents) rather than zeroed, matching the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 6 - refreshed metrics
$ws.Range("D2").Value = 0.00379
$ws.Range("E2").Value = 0.0394
$ws.Range("F2").Value = -0.007800000000000001
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1774
$ws.Range("L2").Value = 0.2064975730133048
$ws.Range("M2").Value = 1669.488
$ws.Range("N2").Value = 0.0880527001439866
$ws.Range("O2").Value = 0.941086809470124
$ws.Range("P2").Value = 168.62
$ws.Range("Q2").Value = 0.008893413009425057
$ws.Range("R2").Value = 0.09505073280721534
$ws.Range("S2").Value = 1500.868
$ws.Range("T2").Value = 0.8989989745359056
$ws.Range("U2").Value = 19065.1
$ws.Range("V2").Value = 1.005537945474971
$ws.Range("W2").Value = 0.07400158995881501
$ws.Range("X2").Value = 0.06262870969906026
$ws.Range("Y2").Value = 0.01137288025975475
$ws.Range("Z2").Value = 0.02946901977922161
$ws.Range("AB2").Value = 0.03233898268077744
$ws.Range("AC2").Value = -0.03233898268077744
$ws.Range("AD2").Value = 295326.9
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 295326.9
$ws.Range("AG2").Value = 276261.8
$ws.Range("AH2").Value = 0.9396726558845896
$ws.Range("AI2").Value = 0.8939758839202524
$ws.Range("AJ2").Value = 0.9357767834974303
$ws.Range("AK2").Value = 0.887482327407721
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

# Row 3: "Sparekassen Sjælland-Fyn A/S (CPSE:SPKSJF)" -> "Lån & Spar Bank A/S (CPSE:LASP)"
$ws.Range("B3").Value = "Lån & Spar Bank A/S (CPSE:LASP)"
$ws.Range("D3").Value = 0.06519999999999999
$ws.Range("E3").Value = 0.108
$ws.Range("K3").Value = 24.9
$ws.Range("L3").Value = 0.1809593023255814
$ws.Range("M3").Value = 22.06
$ws.Range("N3").Value = 0.07940964722822175
$ws.Range("O3").Value = 0.8859437751004018
$ws.Range("P3").Value = 5.46
$ws.Range("Q3").Value = 0.01965442764578833
$ws.Range("R3").Value = 0.219277108433735
$ws.Range("S3").Value = 16.6
$ws.Range("T3").Value = 0.7524932003626473
$ws.Range("U3").Value = 1267.7
$ws.Range("V3").Value = 4.563354931605471
$ws.Range("W3").Value = 0.09067734887108521
$ws.Range("X3").Value = 0.03426060826462669
$ws.Range("Y3").Value = 0.05641674060645852
$ws.Range("Z3").Value = -0.1844504021447721
$ws.Range("AA3").Value = -0
$ws.Range("AB3").Value = 0.03200481174797615
$ws.Range("AC3").Value = -0.03200481174797615
$ws.Range("AD3").Value = 36.8
$ws.Range("AF3").Value = 36.8
$ws.Range("AG3").Value = -1230.9
$ws.Range("AH3").Value = 0.1169739351557533
$ws.Range("AI3").Value = 0.1041607698839513
$ws.Range("AJ3").Value = 1.291469940195153
$ws.Range("AK3").Value = 1.346128608923884

# Row 4: "Danske Bank A/S (CPSE:DANSKE)" -> "P/F BankNordik (CPSE:BNORDIK CSE)"
$ws.Range("B4").Value = "P/F BankNordik (CPSE:BNORDIK CSE)"
$ws.Range("D4").Value = 0.00094
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 27
$ws.Range("L4").Value = 0.2586206896551724
$ws.Range("M4").Value = 1.118
$ws.Range("N4").Value = 0.004717299578059072
$ws.Range("O4").Value = 0.04140740740740741
$ws.Range("P4").Value = 1.06
$ws.Range("Q4").Value = 0.004472573839662448
$ws.Range("R4").Value = 0.03925925925925926
$ws.Range("S4").Value = 0.05800000000000005
$ws.Range("T4").Value = 0.05187835420393564
$ws.Range("U4").Value = 220.6
$ws.Range("V4").Value = 0.9308016877637131
$ws.Range("W4").Value = 0.08353960396039604
$ws.Range("X4").Value = 0.03511720174162086
$ws.Range("Y4").Value = 0.04842240221877518
$ws.Range("Z4").Value = 0.4035562427522226
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.03203017939537029
$ws.Range("AC4").Value = -0.03203017939537029
$ws.Range("AD4").Value = 42.9
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 42.9
$ws.Range("AG4").Value = -177.7
$ws.Range("AH4").Value = 0.1532690246516613
$ws.Range("AI4").Value = 0.1025824964131994
$ws.Range("AJ4").Value = -2.99662731871838
$ws.Range("AK4").Value = -0.8992914979757083
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()

# Row 5: "Jyske Bank A/S (CPSE:JYSK)" -> "Sparekassen Sjælland-Fyn A/S (CPSE:SPKSJF)"
$ws.Range("B5").Value = "Sparekassen Sjælland-Fyn A/S (CPSE:SPKSJF)"
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 35.5
$ws.Range("L5").Value = 0.2156743620899149
$ws.Range("M5").Value = 11.31
$ws.Range("N5").Value = 0.04416243654822334
$ws.Range("O5").Value = 0.3185915492957746
$ws.Range("P5").Value = 8.199999999999999
$ws.Range("Q5").Value = 0.03201874267864115
$ws.Range("R5").Value = 0.2309859154929577
$ws.Range("S5").Value = 3.109999999999999
$ws.Range("T5").Value = 0.2749778956675508
$ws.Range("U5").Value = 468.6
$ws.Range("V5").Value = 1.829754002342835
$ws.Range("W5").Value = 0.07627847013321874
$ws.Range("X5").Value = 0.04743994499157295
$ws.Range("Y5").Value = 0.0288385251416458
$ws.Range("Z5").Value = 0.4455874390904168
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0.03225008334166182
$ws.Range("AC5").Value = -0.03225008334166182
$ws.Range("AD5").Value = 225.2
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 225.2
$ws.Range("AG5").Value = -243.4
$ws.Range("AH5").Value = 0.4678994390193226
$ws.Range("AI5").Value = 0.2997870074547391
$ws.Range("AJ5").Value = -19.16535433070868
$ws.Range("AK5").Value = -0.8612880396319889
$ws.Range("AN5").ClearContents()
$ws.Range("AP5").ClearContents()

# Row 6: "Lån & Spar Bank A/S (CPSE:LASP)" -> "Sydbank A/S (CPSE:SYDB)"
$ws.Range("B6").Value = "Sydbank A/S (CPSE:SYDB)"
$ws.Range("D6").Value = 0.00379
$ws.Range("E6").Value = -0.0617
$ws.Range("K6").Value = 121.1
$ws.Range("L6").Value = 0.1868250539956803
$ws.Range("M6").Value = 251.6
$ws.Range("N6").Value = 0.1926493108728943
$ws.Range("O6").Value = 2.077621800165153
$ws.Range("P6").Value = 1.1
$ws.Range("Q6").Value = 0.0008422664624808577
$ws.Range("R6").Value = 0.009083402146985962
$ws.Range("S6").Value = 250.5
$ws.Range("T6").Value = 0.9956279809220986
$ws.Range("U6").Value = 378.9
$ws.Range("V6").Value = 0.2901225114854518
$ws.Range("W6").Value = 0.07172470978441127
$ws.Range("X6").Value = 0.07781747440654757
$ws.Range("Y6").Value = -0.006092764622136301
$ws.Range("Z6").Value = 0.1080099311815773
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 0.03242788201989306
$ws.Range("AC6").Value = -0.03242788201989306
$ws.Range("AD6").Value = 3396.7
$ws.Range("AF6").Value = 3396.7
$ws.Range("AG6").Value = 3017.8
$ws.Range("AH6").Value = 0.7222871967167798
$ws.Range("AI6").Value = 0.6372436823443334
$ws.Range("AJ6").Value = 0.697950876543781
$ws.Range("AK6").Value = 0.6094841862907461

# Row 7: "P/F BankNordik (CPSE:BNORDIK CSE)" -> "Danske Bank A/S (CPSE:DANSKE)"
$ws.Range("B7").Value = "Danske Bank A/S (CPSE:DANSKE)"
$ws.Range("D7").Value = -0.0257
$ws.Range("E7").Value = 0.0365
$ws.Range("F7").Value = -0.007800000000000001
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1290.4
$ws.Range("L7").Value = 0.2039964588339446
$ws.Range("M7").Value = 1009.9
$ws.Range("N7").Value = 0.071606338852058
$ws.Range("O7").Value = 0.7826255424674519
$ws.Range("P7").Value = 126.5
$ws.Range("Q7").Value = 0.008969404757684263
$ws.Range("R7").Value = 0.09803161810291382
$ws.Range("S7").Value = 883.4
$ws.Range("T7").Value = 0.8747400732745816
$ws.Range("U7").Value = 12855.5
$ws.Range("V7").Value = 0.9115113269755735
$ws.Range("W7").Value = 0.05849660462206588
$ws.Range("X7").Value = 0.3047858711372989
$ws.Range("Y7").Value = -0.246289266515233
$ws.Range("Z7").Value = 0.02882603712623695
$ws.Range("AA7").Value = 0
$ws.Range("AB7").Value = 0.03396779263179799
$ws.Range("AC7").Value = -0.03396779263179799
$ws.Range("AD7").Value = 218084.8
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 218084.8
$ws.Range("AG7").Value = 205229.3
$ws.Range("AH7").Value = 0.9392583519496891
$ws.Range("AI7").Value = 0.8926706671141929
$ws.Range("AJ7").Value = 0.9356981719104485
$ws.Range("AK7").Value = 0.8867092531664438
$ws.Range("AN7").ClearContents()
$ws.Range("AP7").ClearContents()

# Row 8: "Sydbank A/S (CPSE:SYDB)" -> "Jyske Bank A/S (CPSE:JYSK)"
$ws.Range("B8").Value = "Jyske Bank A/S (CPSE:JYSK)"
$ws.Range("D8").Value = 0.006079999999999999
$ws.Range("E8").Value = 0.0423
$ws.Range("K8").Value = 275.1
$ws.Range("L8").Value = 0.2272614622057001
$ws.Range("M8").Value = 373.5
$ws.Range("N8").Value = 0.1343670180235277
$ws.Range("O8").Value = 1.357688113413304
$ws.Range("P8").Value = 26.3
$ws.Range("Q8").Value = 0.009461452674749075
$ws.Range("R8").Value = 0.0956015994183933
$ws.Range("S8").Value = 347.2
$ws.Range("T8").Value = 0.9295850066934404
$ws.Range("U8").Value = 3873.8
$ws.Range("V8").Value = 1.393603626290607
$ws.Range("W8").Value = 0.05349434138373586
$ws.Range("X8").Value = 0.4987699917984738
$ws.Range("Y8").Value = -0.4452756504147379
$ws.Range("Z8").Value = 0.01828571946488779
$ws.Range("AB8").Value = 0.03402073675419914
$ws.Range("AC8").Value = -0.03402073675419914
$ws.Range("AD8").Value = 73540.5
$ws.Range("AF8").Value = 73540.5
$ws.Range("AG8").Value = 69666.7
$ws.Range("AH8").Value = 0.963578449742008
$ws.Range("AI8").Value = 0.9286213967865928
$ws.Range("AJ8").Value = 0.9616309437045871
$ws.Range("AK8").Value = 0.9249502784143262
